# The sheet holds one daily price record per row (rows 249-351). A new
# daily record was inserted at row 249 (date serial 44795 = 2022-08-22),
# pushing every existing record from row 249 onward down by one row - the
# former row 351 record now lives at row 352. This grows the used range
# from A1:R351 to A1:R352.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 249..351 down to 250..352, carrying values/formatting along.
$ws.Rows.Item(249).Insert()

# Populate the newly opened row 249 with the new record.
$ws.Cells.Item(249, 1).Value = 10
$ws.Cells.Item(249, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(249, 3).Value = "La Araucanía"
$ws.Cells.Item(249, 4).Value = 44795
$ws.Cells.Item(249, 5).Value = 9
$ws.Cells.Item(249, 6).Value = 100112017
$ws.Cells.Item(249, 7).Value = "Apio"
$ws.Cells.Item(249, 8).Value = "Americana (o)"
$ws.Cells.Item(249, 9).Value = "Primera"
$ws.Cells.Item(249, 10).Value = 380
$ws.Cells.Item(249, 11).Value = 12000
$ws.Cells.Item(249, 12).Value = 13000
$ws.Cells.Item(249, 13).Value = 12474
$ws.Cells.Item(249, 14).Value = "`$/docena de matas"
$ws.Cells.Item(249, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(249, 16).Value = 2079
$ws.Cells.Item(249, 17).Value = 6
$ws.Cells.Item(249, 18).Value = "Hortaliza"
